# Invalid-credential login test: add a new "BookStoreInvalidLogin" sheet
# (with test data for the BookStore login negative-test cases) right
# after "BookStoreLogin" and before "Alerts".

$wb = $excel.ActiveWorkbook

$bookStoreLogin = $wb.Worksheets.Item("BookStoreLogin")

# New sheet is inserted immediately after BookStoreLogin (pushes Alerts
# one slot to the right).
$newSheet = $wb.Worksheets.Add($null, $bookStoreLogin)
$newSheet.Name = "BookStoreInvalidLogin"

# Header row.
$newSheet.Range("A1").Value = "Invalid username"
$newSheet.Range("B1").Value = "Invalid password"
$newSheet.Range("C1").Value = "Error message"

# Test-case rows.
$newSheet.Range("A2").Value = "johan123"
$newSheet.Range("B2").Value = "jedandva123!"
$newSheet.Range("C2").Value = "Invalid username or password!"

$newSheet.Range("A3").Value = "JOHAN1235"
$newSheet.Range("B3").Value = "Jedandva12!"

$newSheet.Range("A4").Value = "johan1235"
$newSheet.Range("B4").Value = "JEDANdva12!"

$newSheet.Range("A5").Value = "johan123"
$newSheet.Range("B5").Value = "Jedan1dva2!"

# Header styling: Excel's built-in "Neutral" cell style (orange text on
# pale-yellow fill).
$newSheet.Range("A1:C1").Style = "Neutral"

# Data cells are stored as explicit Text format, matching the rest of the
# workbook's test-data sheets.
$newSheet.Range("A2:B5").NumberFormat = "@"

# Column widths matching the authored sheet.
$newSheet.Columns.Item(1).ColumnWidth = 21.833333333333332
$newSheet.Columns.Item(2).ColumnWidth = 21.5

# BookStoreLogin's selection/active-cell moves off the old C2 spot.
[void]$bookStoreLogin.Activate()
$bookStoreLogin.Range("A2").Select() | Out-Null

# The newly-added sheet becomes the active tab, selection on A4.
[void]$newSheet.Activate()
$newSheet.Range("A4").Select() | Out-Null
